$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C; existing column C (dates/ratings)
# shifts to column E.
$ws.Range("C:D").Insert()

# Row 1 date headers: newest dates added to the left (B, C), older ones
# shift right (D, E).
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = "Jun_13"
$ws.Range("E1").Value = "Jun_10"

# Fill rows 2-27 with the "UN" rating in the new C/D columns, and the
# shifted E column, matching column B.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
    $ws.Cells.Item($r, 5).Value = "UN"
}

# Cosmetic column widths (matches the pre-existing column C width).
$ws.Columns.Item(3).ColumnWidth = 7.1666
$ws.Columns.Item(4).ColumnWidth = 7.1666
$ws.Columns.Item(5).ColumnWidth = 7.1666
